$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 140
$ws.Cells.Item(2, 9).Value = 63.42857
$ws.Cells.Item(2, 11).Value = 63.42857
$ws.Cells.Item(2, 13).Value = 49.57143
$ws.Cells.Item(98, 8).Value = 1568.7838
$ws.Cells.Item(98, 9).Value = 1838.1333
$ws.Cells.Item(98, 11).Value = 1838.1333
$ws.Cells.Item(98, 13).Value = -340.1333
$ws.Cells.Item(122, 8).Value = 1568.7838
$ws.Cells.Item(122, 9).Value = 1838.1333
$ws.Cells.Item(122, 11).Value = 5514.3999
$ws.Cells.Item(122, 13).Value = -3064.3999

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3087.3076
$ws.Cells.Item(61, 9).Value = 1813.1578
$ws.Cells.Item(61, 11).Value = 1813.1578
$ws.Cells.Item(61, 13).Value = -1601.1578
$ws.Cells.Item(74, 8).Value = 27780570
$ws.Cells.Item(74, 9).Value = 37038396
$ws.Cells.Item(74, 11).Value = 37038396
$ws.Cells.Item(74, 13).Value = -37037522
$ws.Cells.Item(77, 8).Value = 27780570
$ws.Cells.Item(77, 9).Value = 37038396
$ws.Cells.Item(77, 11).Value = 185191980
$ws.Cells.Item(77, 13).Value = -185187612
$ws.Cells.Item(93, 8).Value = 30000
$ws.Cells.Item(93, 10).Value = 30000
$ws.Cells.Item(93, 12).Value = 30000
$ws.Cells.Item(93, 14).Value = -33744
$ws.Cells.Item(102, 8).Value = 2076197.6
$ws.Cells.Item(102, 9).Value = 2639578.5
$ws.Cells.Item(102, 10).Value = 104364.3
$ws.Cells.Item(102, 11).Value = 2639578.5
$ws.Cells.Item(102, 12).Value = 104364.3
$ws.Cells.Item(102, 13).Value = -2637956.5
$ws.Cells.Item(102, 14).Value = -107608.3
$ws.Cells.Item(132, 8).Value = 5858.851
$ws.Cells.Item(132, 9).Value = 5458.8647
$ws.Cells.Item(132, 11).Value = 16376.5941
$ws.Cells.Item(132, 13).Value = -13846.5941
$ws.Cells.Item(136, 8).Value = 3087.3076
$ws.Cells.Item(136, 9).Value = 1813.1578
$ws.Cells.Item(136, 11).Value = 5439.4734
$ws.Cells.Item(136, 13).Value = -2889.4734
$ws.Cells.Item(138, 8).Value = 70939.75
$ws.Cells.Item(138, 9).Value = 48769
$ws.Cells.Item(138, 11).Value = 48769
$ws.Cells.Item(138, 13).Value = -43629

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 4236.1304
$ws.Cells.Item(134, 9).Value = 3542.5881
$ws.Cells.Item(134, 11).Value = 10627.7643
$ws.Cells.Item(134, 13).Value = -8092.764299999999
$ws.Cells.Item(136, 8).Value = 3000
$ws.Cells.Item(136, 9).Value = 3000
$ws.Cells.Item(136, 11).Value = 9000
$ws.Cells.Item(136, 13).Value = -3900

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1284.3572
$ws.Cells.Item(16, 9).Value = 1128.8
$ws.Cells.Item(16, 11).Value = 1128.8
$ws.Cells.Item(16, 13).Value = -841.8
$ws.Cells.Item(31, 8).Value = 3429.4307
$ws.Cells.Item(31, 10).Value = 3491.6492
$ws.Cells.Item(31, 12).Value = 3491.6492
$ws.Cells.Item(31, 14).Value = -4081.6492
$ws.Cells.Item(34, 8).Value = 3429.4307
$ws.Cells.Item(34, 10).Value = 3491.6492
$ws.Cells.Item(34, 12).Value = 3491.6492
$ws.Cells.Item(34, 14).Value = -3895.6492
$ws.Cells.Item(58, 8).Value = 2034.2858
$ws.Cells.Item(58, 9).Value = 1299.3125
$ws.Cells.Item(58, 11).Value = 1299.3125
$ws.Cells.Item(58, 13).Value = -1096.3125
$ws.Cells.Item(113, 8).Value = 1284.3572
$ws.Cells.Item(113, 9).Value = 1128.8
$ws.Cells.Item(113, 11).Value = 1128.8
$ws.Cells.Item(113, 13).Value = 1041.2
$ws.Cells.Item(134, 8).Value = 2458.311
$ws.Cells.Item(134, 9).Value = 1164.3243
$ws.Cells.Item(134, 10).Value = 8443
$ws.Cells.Item(134, 11).Value = 3492.9729
$ws.Cells.Item(134, 12).Value = 25329
$ws.Cells.Item(134, 13).Value = -957.9728999999998
$ws.Cells.Item(134, 14).Value = -30399
$ws.Cells.Item(136, 8).Value = 2034.2858
$ws.Cells.Item(136, 9).Value = 1299.3125
$ws.Cells.Item(136, 11).Value = 3897.9375
$ws.Cells.Item(136, 13).Value = -1347.9375

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 148
$ws.Cells.Item(12, 10).Value = 173.13333
$ws.Cells.Item(12, 12).Value = 519.39999
$ws.Cells.Item(12, 14).Value = -865.39999
$ws.Cells.Item(23, 8).Value = 2388.7407
$ws.Cells.Item(23, 10).Value = 2880.762
$ws.Cells.Item(23, 12).Value = 8642.286
$ws.Cells.Item(23, 14).Value = -9112.286
$ws.Cells.Item(122, 8).Value = 66666868
$ws.Cells.Item(122, 10).Value = 200000000
$ws.Cells.Item(122, 12).Value = 1800000000
$ws.Cells.Item(122, 14).Value = -1800004900
$ws.Cells.Item(132, 8).Value = 2810.087
$ws.Cells.Item(132, 9).Value = 1950
$ws.Cells.Item(132, 10).Value = 2955.8645
$ws.Cells.Item(132, 11).Value = 17550
$ws.Cells.Item(132, 12).Value = 26602.7805
$ws.Cells.Item(132, 13).Value = -15020
$ws.Cells.Item(132, 14).Value = -31662.7805

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 710.3333
$ws.Cells.Item(97, 9).Value = 662.76
$ws.Cells.Item(97, 10).Value = 948.2
$ws.Cells.Item(97, 11).Value = 662.76
$ws.Cells.Item(97, 12).Value = 948.2
$ws.Cells.Item(97, 13).Value = -166.76
$ws.Cells.Item(97, 14).Value = -1940.2

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3691
$ws.Cells.Item(7, 9).Value = 3299.875
$ws.Cells.Item(7, 11).Value = 3299.875
$ws.Cells.Item(7, 13).Value = -3187.875
$ws.Cells.Item(61, 8).Value = 4510.1
$ws.Cells.Item(61, 9).Value = 3762
$ws.Cells.Item(61, 11).Value = 3762
$ws.Cells.Item(61, 13).Value = -3560
$ws.Cells.Item(82, 8).Value = 744
$ws.Cells.Item(82, 9).Value = 715.2
$ws.Cells.Item(82, 11).Value = 715.2
$ws.Cells.Item(82, 13).Value = -354.2
$ws.Cells.Item(85, 8).Value = 744
$ws.Cells.Item(85, 9).Value = 715.2
$ws.Cells.Item(85, 11).Value = 715.2
$ws.Cells.Item(85, 13).Value = 532.8
$ws.Cells.Item(93, 8).Value = 83335336
$ws.Cells.Item(93, 10).Value = 2752
$ws.Cells.Item(93, 12).Value = 2752
$ws.Cells.Item(93, 14).Value = -5248
$ws.Cells.Item(96, 8).Value = 74241
$ws.Cells.Item(96, 9).Value = 53222
$ws.Cells.Item(96, 10).Value = 84750.5
$ws.Cells.Item(96, 11).Value = 53222
$ws.Cells.Item(96, 12).Value = 84750.5
$ws.Cells.Item(96, 13).Value = -50476
$ws.Cells.Item(96, 14).Value = -90242.5
$ws.Cells.Item(113, 8).Value = 4510.1
$ws.Cells.Item(113, 9).Value = 3762
$ws.Cells.Item(113, 11).Value = 3762
$ws.Cells.Item(113, 13).Value = -1592
$ws.Cells.Item(126, 8).Value = 3691
$ws.Cells.Item(126, 9).Value = 3299.875
$ws.Cells.Item(126, 11).Value = 9899.625
$ws.Cells.Item(126, 13).Value = -7429.625
$ws.Cells.Item(132, 8).Value = 43481690
$ws.Cells.Item(132, 10).Value = 4142.636
$ws.Cells.Item(132, 12).Value = 12427.908
$ws.Cells.Item(132, 14).Value = -17487.908
$ws.Cells.Item(134, 8).Value = 0
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 14).ClearContents()
$ws.Cells.Item(135, 8).Value = 47143
$ws.Cells.Item(135, 10).Value = 47143
$ws.Cells.Item(135, 12).Value = 47143
$ws.Cells.Item(135, 14).Value = -57283
$ws.Cells.Item(136, 8).Value = 4051.0156
$ws.Cells.Item(136, 9).Value = 3489.3914
$ws.Cells.Item(136, 10).Value = 5486.278
$ws.Cells.Item(136, 11).Value = 10468.1742
$ws.Cells.Item(136, 12).Value = 16458.834
$ws.Cells.Item(136, 13).Value = -7918.174199999999
$ws.Cells.Item(136, 14).Value = -21558.834
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 14).ClearContents()
$ws.Cells.Item(139, 8).Value = 58801.6
$ws.Cells.Item(139, 10).Value = 58801.6
$ws.Cells.Item(139, 12).Value = 58801.6
$ws.Cells.Item(139, 14).Value = -69081.6

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(123, 8).Value = 42232.668
$ws.Cells.Item(123, 9).Value = 38699
$ws.Cells.Item(123, 10).Value = 43999.5
$ws.Cells.Item(123, 11).Value = 38699
$ws.Cells.Item(123, 12).Value = 43999.5
$ws.Cells.Item(123, 13).Value = -33799
$ws.Cells.Item(123, 14).Value = -53799.5
$ws.Cells.Item(132, 8).Value = 7034.9585
$ws.Cells.Item(132, 9).Value = 8763.333
$ws.Cells.Item(132, 11).Value = 26289.999
$ws.Cells.Item(132, 13).Value = -23759.999
$ws.Cells.Item(136, 8).Value = 11241.286
$ws.Cells.Item(136, 9).Value = 15801.375
$ws.Cells.Item(136, 11).Value = 47404.125
$ws.Cells.Item(136, 13).Value = -44854.125
